$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5482547283172607
$ws.Range("B1").Value = 1.760244727134705
$ws.Range("C1").Value = 4.940728187561035
$ws.Range("D1").Value = 1.690102458000183
$ws.Range("E1").Value = 0.8968498706817627
